$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.863.33"
$ws.Range("E2").Value = "  -2.12%  "

$ws.Range("D3").Value = "3.503.42"
$ws.Range("E3").Value = "  -3.59%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.11%  "

$ws.Range("D7").Value = "3.501.91"
$ws.Range("E7").Value = "  -3.59%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.488"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.144"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.92%  "

$ws.Range("E11").Value = "  +4.49%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.430"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.31%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000215"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.29%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.32%  "

$ws.Range("D15").Value = "4.085.91"
$ws.Range("E15").Value = "  -4.01%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "67.748.25"
$ws.Range("E16").Value = "  -2.20%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.489.49"
$ws.Range("E17").Value = "  -4.35%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.117"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.35%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.77%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.34%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.68%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "444.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.36%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.626"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.56%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.54%  "

$ws.Range("D25").Value = "3.633.94"
$ws.Range("E25").Value = "  -4.04%  "

$ws.Range("E26").Value = "  -0.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000124"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -9.29%  "

$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.97"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.26%  "

$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.67"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.32%  "

$ws.Range("E31").Value = "  -3.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.170"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.86%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.59"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.23%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.17"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.08%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.85"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.27%  "

$ws.Range("B37").Value = "RenzoRestakedETH"
$ws.Range("C37").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D37").Value = "3.488.96"
$ws.Range("E37").Value = "  -3.83%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.61%  "

$ws.Range("E39").Value = "  +0.00%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.56%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.998"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.21%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "174.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.26%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0903"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.73%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.25%  "

$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "30.60"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.67%  "

$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.898"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.39%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.62"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.32%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -11.91%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.62%  "
